# Carjacking arrests by month YoY - add data point for 2021-10-10
# (update "through" date label from 10-01 to 10-02 and refresh October + Total figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet / workbook tab name: "Through 2021-10-01" -> "Through 2021-10-02"
$ws.Name = "Through 2021-10-02"

# --- Row 12 ("October (through 10-01)" -> "October (through 10-02)") ---
$ws.Range("A12").Value = "October (through 10-02)"

$ws.Range("C12").Value = 3
$ws.Range("F12").Value = 6

$ws.Range("H12").Value = 2
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 0.25
$ws.Range("L12").Value = 5
$ws.Range("O12").Value = 1
$ws.Range("R12").Value = 9
$ws.Range("U12").Value = 17

# --- Row 13 (Total row) ---
$ws.Range("C13").Value = 199
$ws.Range("D13").Value = 0.131

$ws.Range("F13").Value = 389
$ws.Range("G13").Value = 0.1057

$ws.Range("H13").Value = 52
$ws.Range("I13").Value = 583
$ws.Range("J13").Value = 0.0819

$ws.Range("L13").Value = 492
$ws.Range("M13").Value = 0.1103

$ws.Range("R13").Value = 857
$ws.Range("S13").Value = 0.0582

$ws.Range("U13").Value = 1187
$ws.Range("V13").Value = 0.0617
